# Apply "sliding window" update to the dataset:
#  - drop the oldest 8 data rows (old rows 2-9)
#  - keep the remaining 13 data rows (old rows 10-22), shifted up to rows 2-14
#  - append 7 brand new data rows (new rows 15-21)
#  - the sheet now spans A1:C21 (was A1:C22)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full dataset (20 rows) to place under the header in A2:C21
$data = @(
    @(-0.1843285858631134, -2.109012365341187, -1.700190663337708),
    @(-0.3310889601707458, -3.295770645141602, -1.143387079238892),
    @(-0.3715587854385376, -2.42528772354126, -1.059240341186523),
    @(-0.6349944472312927, -7.78822660446167, -1.318247199058533),
    @(0.4980078935623169, -0.0974330082535743, 0.2806925773620605),
    @(-0.831234872341156, 2.706896543502808, 1.08245325088501),
    @(0.7096726894378662, 3.787517309188842, 0.2267837226390838),
    @(-1.536326050758362, 1.956906795501709, 0.3100140988826751),
    @(0.8593347072601318, 7.671703815460205, 1.616807579994202),
    @(5.033376693725586, 4.168544292449951, 0.3526219725608825),
    @(-0.6580545902252197, 1.885893702507019, 0.1637118905782699),
    @(-2.827084302902222, -3.540880441665649, -0.5829181671142578),
    @(-1.010065674781799, -4.767192363739014, -0.09926560521125791),
    @(-1.965306162834168, -0.784503698348999, 1.217912554740906),
    @(-0.3843869566917419, -3.860210180282593, 2.725528001785278),
    @(4.380514621734619, -3.570049285888672, 1.001513600349426),
    @(0.96409809589386, -0.7756461501121521, -0.2683225572109222),
    @(-0.6479753255844116, 0.6890559792518616, 1.091005325317383),
    @(-1.356426239013672, 3.433979034423828, -1.384373307228088),
    @(0.2755002379417419, 2.776687860488892, -1.657124638557434)
)

$rowCount = $data.Length
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}

# Remove the now-stale trailing row (old row 22), since the data shrank by one row overall.
$ws.Range("A22:C22").Delete()
